$d = $word.ActiveDocument

function New-WordPkgXml([string]$bodyXml) {
    return "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" +
           "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" +
           "<pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
           "<w:body>$bodyXml</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
}

# 1) First paragraph: drop the Titre/lev paragraph formatting entirely and
#    split the new text across two runs (as the edit introduced two <w:r>).
$p1 = $d.Paragraphs(1)
$p1.Range.InsertXML((New-WordPkgXml "<w:p><w:r><w:t>MAIS NIQUE TA Me</w:t></w:r><w:r><w:t>RE WORD</w:t></w:r></w:p>"))

# 2) Second paragraph text swap.
$d.Content.Find.Execute("J'aime trop les frites putain", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "NBGFJDLGBLDFSPOGS", 2) | Out-Null

# 3) Third paragraph text swap.
$d.Content.Find.Execute("Cfds", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Pourquoi", 2) | Out-Null

# 4) Fourth paragraph text swap.
$d.Content.Find.Execute("Fds", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Pourauoi", 2) | Out-Null

# 5) Remove the "Gh" / "Hj" / "Kj" paragraphs entirely (including their marks).
$d.Paragraphs(5).Range.Delete() | Out-Null
$d.Paragraphs(5).Range.Delete() | Out-Null
$d.Paragraphs(5).Range.Delete() | Out-Null

# 6) Last paragraph ("hdgsfad" -> "p", bookmark kept) plus a brand-new empty
#    trailing paragraph, written together so the empty one stays truly empty.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertXML((New-WordPkgXml ("<w:p><w:r><w:t>p</w:t></w:r>" +
    "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p><w:p/>")))
